$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.914.96"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "3.127.25"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'589.58"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").Value = "'135.78"
$ws.Range("E6").Value = "  -4.94%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.120.92"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("E10").Value = "  -4.12%  "
$ws.Range("D11").Value = "'5.24"
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D12").Value = "'0.453"
$ws.Range("E12").Value = "  -3.18%  "
$ws.Range("D13").Value = "'0.0000242"
$ws.Range("E13").Value = "  -5.42%  "
$ws.Range("D14").Value = "'33.90"
$ws.Range("E14").Value = "  -3.54%  "
$ws.Range("D15").Value = "3.635.23"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").Value = "62.951.62"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").Value = "3.118.65"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").Value = "  -4.44%  "
$ws.Range("D20").Value = "'468.52"
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("D21").Value = "'14.04"
$ws.Range("E21").Value = "  -3.53%  "
$ws.Range("D22").Value = "'0.694"
$ws.Range("E22").Value = "  -2.19%  "
$ws.Range("D23").Value = "'7.64"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("D24").Value = "'85.12"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Value = "'12.89"
$ws.Range("E25").Value = "  -3.94%  "
$ws.Range("D27").Value = "'2.70"
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").Value = "'7.84"
$ws.Range("E28").Value = "  -6.11%  "
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("D30").Value = "'6.78"
$ws.Range("E30").Value = "  -4.98%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'26.51"
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("D33").Value = "'0.108"
$ws.Range("E33").Value = "  -4.37%  "
$ws.Range("E34").Value = "  -4.40%  "
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("D36").Value = "'51.94"
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").Value = "'5.73"
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("D38").Value = "0.0₃0676"
$ws.Range("E38").Value = "  -12.25%  "
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("D40").Value = "'414.00"
$ws.Range("E40").Value = "  -7.04%  "
$ws.Range("D41").Value = "'8.16"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("D42").Value = "2.901.71"
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("D43").Value = "'2.66"
$ws.Range("E43").Value = "  -11.70%  "
$ws.Range("E44").Value = "  -6.58%  "
$ws.Range("D45").Value = "'0.256"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'2.09"
$ws.Range("E47").Value = "  -5.82%  "
$ws.Range("D48").Value = "'25.26"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("E50").Value = "  -8.64%  "
$ws.Range("D51").Value = "'120.42"
$ws.Range("E51").Value = "  -0.10%  "
